$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.195046439628483
$ws.Range("C2").Value = 0.544891640866873
$ws.Range("J2").Value = 0.01238390092879257
$ws.Range("O2").Value = 0.003095975232198143
$ws.Range("P2").Value = 0.1578947368421053
$ws.Range("S2").Value = 0.08668730650154799
$ws.Range("B3").Value = 0.02185792349726776
$ws.Range("C3").Value = 0.03825136612021858
$ws.Range("J3").Value = 0.02185792349726776
$ws.Range("P3").Value = 0.7431693989071039
$ws.Range("S3").Value = 0.1748633879781421
$ws.Range("J4").Value = 0.04347826086956522
$ws.Range("P4").Value = 0.6521739130434783
$ws.Range("S4").Value = 0.3043478260869565
$ws.Range("B6").Value = 0.05508474576271186
$ws.Range("D6").Value = 0.01271186440677966
$ws.Range("F6").Value = 0.05932203389830509
$ws.Range("J6").Value = 0.288135593220339
$ws.Range("O6").Value = 0.0211864406779661
$ws.Range("Q6").Value = 0.1525423728813559
$ws.Range("R6").Value = 0.07627118644067797
$ws.Range("S6").Value = 0.3347457627118644
$ws.Range("B7").Value = 0.1036585365853658
$ws.Range("D7").Value = 0.01219512195121951
$ws.Range("E7").Value = 0.006097560975609756
$ws.Range("F7").Value = 0.05487804878048781
$ws.Range("J7").Value = 0.1524390243902439
$ws.Range("O7").Value = 0.02439024390243903
$ws.Range("Q7").Value = 0.1768292682926829
$ws.Range("R7").Value = 0.04878048780487805
$ws.Range("S7").Value = 0.4207317073170732
$ws.Range("B8").Value = 0.1044083526682135
$ws.Range("D8").Value = 0.01160092807424594
$ws.Range("E8").Value = 0.002320185614849188
$ws.Range("F8").Value = 0.0765661252900232
$ws.Range("J8").Value = 0.09280742459396751
$ws.Range("O8").Value = 0.03016241299303944
$ws.Range("Q8").Value = 0.2018561484918794
$ws.Range("R8").Value = 0.0951276102088167
$ws.Range("S8").Value = 0.3851508120649652
$ws.Range("B9").Value = 0.1
$ws.Range("D9").Value = 0.02222222222222222
$ws.Range("F9").Value = 0.06666666666666667
$ws.Range("J9").Value = 0.08333333333333333
$ws.Range("O9").Value = 0.005555555555555556
$ws.Range("Q9").Value = 0.2444444444444444
$ws.Range("R9").Value = 0.05
$ws.Range("S9").Value = 0.4277777777777778
$ws.Range("B10").Value = 0.1184905660377358
$ws.Range("D10").Value = 0.02490566037735849
$ws.Range("E10").Value = 0.002264150943396227
$ws.Range("F10").Value = 0.06566037735849056
$ws.Range("J10").Value = 0.1184905660377358
$ws.Range("O10").Value = 0.01509433962264151
$ws.Range("Q10").Value = 0.1909433962264151
$ws.Range("R10").Value = 0.1049056603773585
$ws.Range("S10").Value = 0.3592452830188679
$ws.Range("G11").Value = 0.1538461538461539
$ws.Range("J11").Value = 0.08461538461538462
$ws.Range("K11").Value = 0.2076923076923077
$ws.Range("L11").Value = 0.5423076923076923
$ws.Range("S11").Value = 0.01153846153846154
$ws.Range("G12").Value = 0.7027027027027027
$ws.Range("J12").Value = 0.2162162162162162
$ws.Range("K12").Value = 0.006756756756756757
$ws.Range("L12").Value = 0.0472972972972973
$ws.Range("S12").Value = 0.02702702702702703
$ws.Range("G13").Value = 0.6388888888888888
$ws.Range("J13").Value = 0.2777777777777778
$ws.Range("S13").Value = 0.08333333333333333
$ws.Range("F15").Value = 0.02371541501976284
$ws.Range("H15").Value = 0.1541501976284585
$ws.Range("I15").Value = 0.05138339920948617
$ws.Range("J15").Value = 0.4071146245059288
$ws.Range("K15").Value = 0.05533596837944664
$ws.Range("M15").Value = 0.01185770750988142
$ws.Range("O15").Value = 0.07114624505928854
$ws.Range("S15").Value = 0.225296442687747
$ws.Range("F16").Value = 0.02830188679245283
$ws.Range("H16").Value = 0.1179245283018868
$ws.Range("I16").Value = 0.04716981132075472
$ws.Range("J16").Value = 0.4575471698113208
$ws.Range("K16").Value = 0.1084905660377359
$ws.Range("M16").Value = 0.01415094339622642
$ws.Range("O16").Value = 0.08018867924528301
$ws.Range("S16").Value = 0.1462264150943396
$ws.Range("F17").Value = 0.0267260579064588
$ws.Range("H17").Value = 0.1826280623608018
$ws.Range("I17").Value = 0.08908685968819599
$ws.Range("J17").Value = 0.4409799554565701
$ws.Range("K17").Value = 0.0957683741648107
$ws.Range("M17").Value = 0.0133630289532294
$ws.Range("N17").Value = 0.0022271714922049
$ws.Range("O17").Value = 0.07572383073496659
$ws.Range("S17").Value = 0.07349665924276169
$ws.Range("F18").Value = 0.05188679245283019
$ws.Range("H18").Value = 0.160377358490566
$ws.Range("I18").Value = 0.09433962264150944
$ws.Range("J18").Value = 0.4056603773584906
$ws.Range("K18").Value = 0.08018867924528301
$ws.Range("M18").Value = 0.02358490566037736
$ws.Range("O18").Value = 0.07547169811320754
$ws.Range("S18").Value = 0.1084905660377359
$ws.Range("F19").Value = 0.01918265221017515
$ws.Range("H19").Value = 0.2093411175979983
$ws.Range("I19").Value = 0.08173477898248541
$ws.Range("J19").Value = 0.4003336113427857
$ws.Range("K19").Value = 0.08924103419516263
$ws.Range("M19").Value = 0.0158465387823186
$ws.Range("O19").Value = 0.07673060884070058
$ws.Range("S19").Value = 0.1075896580483736
